$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report date range (From/To) in row 1
$ws.Range("I1").Value = "07-10-2025 00:00:00"
$ws.Range("K1").Value = "07-10-2025 00:00:00"

# Row 154
$ws.Range("B154").Value = 53925
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44

# Row 155
$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44

# Row 156
$ws.Range("B156").Value = 57756
$ws.Range("E156").Value = 79.37
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644

# Row 256
$ws.Range("B256").Value = 48719
$ws.Range("E256").Value = 353.35
$ws.Range("F256").Value = -81
$ws.Range("G256").Value = -23955.75

# Row 257
$ws.Range("B257").Value = 64979
$ws.Range("E257").Value = 314.41
$ws.Range("F257").Value = 82
$ws.Range("G257").Value = 24251.5

# Row 308
$ws.Range("B308").Value = 63565
$ws.Range("D308").Value = 102.71
$ws.Range("E308").Value = 109.19
$ws.Range("F308").Value = 60
$ws.Range("G308").Value = 6162.6

# Row 309
$ws.Range("B309").Value = 57077
$ws.Range("D309").Value = 93.08
$ws.Range("E309").Value = 111.2
$ws.Range("F309").Value = 1
$ws.Range("G309").Value = 93.08

# Row 310
$ws.Range("B310").Value = 61610
$ws.Range("E310").Value = 122.71
$ws.Range("F310").Value = -58
$ws.Range("G310").Value = -5957.18

# Row 343
$ws.Range("B343").Value = 63531
$ws.Range("F343").Value = 80
$ws.Range("G343").Value = 11478.4

# Row 344
$ws.Range("B344").Value = 63571
$ws.Range("F344").Value = 29
$ws.Range("G344").Value = 4160.92

# Row 347
$ws.Range("B347").Value = 63510
$ws.Range("E347").Value = 50.66
$ws.Range("F347").Value = 167
$ws.Range("G347").Value = 7955.88

# Row 348
$ws.Range("B348").Value = 55356
$ws.Range("E348").Value = 54.04
$ws.Range("F348").Value = -158
$ws.Range("G348").Value = -7527.12

# Row 367
$ws.Range("B367").Value = 61605
$ws.Range("E367").Value = 133.78
$ws.Range("F367").Value = -13
$ws.Range("G367").Value = -1455.48

# Row 368
$ws.Range("B368").Value = 63563
$ws.Range("E368").Value = 119.04
$ws.Range("F368").Value = 15
$ws.Range("G368").Value = 1679.4

# Row 371
$ws.Range("B371").Value = 63564
$ws.Range("E371").Value = 137.16
$ws.Range("F371").Value = 57
$ws.Range("G371").Value = 7353.57

# Row 372
$ws.Range("B372").Value = 61608
$ws.Range("E372").Value = 154.12
$ws.Range("F372").Value = -56
$ws.Range("G372").Value = -7224.56

# Row 528
$ws.Range("B528").Value = 47097
$ws.Range("D528").Value = 112.28
$ws.Range("E528").Value = 134.16
$ws.Range("F528").Value = 15
$ws.Range("G528").Value = 1684.2

# Row 529
$ws.Range("B529").Value = 58047
$ws.Range("D529").Value = 105.54
$ws.Range("E529").Value = 126.1
$ws.Range("F529").Value = 54
$ws.Range("G529").Value = 5699.16

# Row 578
$ws.Range("B578").Value = 64915
$ws.Range("E578").Value = 20.98
$ws.Range("F578").Value = 40
$ws.Range("G578").Value = 789.2

# Row 579
$ws.Range("B579").Value = 45695
$ws.Range("E579").Value = 23.58
$ws.Range("F579").Value = -36
$ws.Range("G579").Value = -710.28

# Row 585
$ws.Range("B585").Value = 45718
$ws.Range("E585").Value = 19.38
$ws.Range("F585").Value = -294
$ws.Range("G585").Value = -4768.68

# Row 586
$ws.Range("B586").Value = 64927
$ws.Range("E586").Value = 17.26
$ws.Range("F586").Value = 295
$ws.Range("G586").Value = 4784.9

# Row 593
$ws.Range("B593").Value = 64919
$ws.Range("E593").Value = 27.97
$ws.Range("F593").Value = 224
$ws.Range("G593").Value = 5891.2

# Row 594
$ws.Range("B594").Value = 45702
$ws.Range("E594").Value = 31.43
$ws.Range("F594").Value = -215
$ws.Range("G594").Value = -5654.5

# Row 679
$ws.Range("B679").Value = 64810
$ws.Range("E679").Value = 291.22
$ws.Range("F679").Value = 7
$ws.Range("G679").Value = 1917.44

# Row 680
$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52

# Row 701
$ws.Range("B701").Value = 60025
$ws.Range("E701").Value = 37.22
$ws.Range("F701").Value = -98
$ws.Range("G701").Value = -3217.34

# Row 702
$ws.Range("B702").Value = 64833
$ws.Range("E702").Value = 34.9
$ws.Range("F702").Value = 99
$ws.Range("G702").Value = 3250.17

# Row 712
$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79

# Row 713
$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11

# Row 864
$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53

# Row 865
$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27
